$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# "bug with empty notes is resolved"
#
# A couple of footnotes in this document ended up with no real content -
# their body is just the bare shad "།" left over from the source
# markup. One of those placeholders (the very last footnote) is pure noise
# and gets removed outright; the other one sits in the middle of the text
# and is filled in with its real content. While we're at it, a stray
# trailing "a" left over in another note is cleaned up too.
# ---------------------------------------------------------------------------

$emptyMarker = "།"
$fn25Text = "བདུད་སྡེ། ཞེས་པར་མ་གཞན་ནང་མེད།"
$fn29Text = "སྒྱུར། སྣར་ཐང་། པེ་ཅིན།"

$count = $d.Footnotes.Count

for ($i = $count; $i -ge 1; $i--) {
    $fn = $d.Footnotes.Item($i)
    $txt = $fn.Range.Text

    if ($txt -eq $emptyMarker) {
        if ($i -eq $count) {
            # Trailing empty placeholder note at the very end of the text -
            # nothing worth keeping, drop the reference and the note.
            $fn.Delete()
        } else {
            # Empty placeholder note earlier in the text - it belongs right
            # after "...དུས་ཀུན་དུ་ནི་བདུད་སྡེ་..."; fill in its real content.
            $fn.Range.Text = $fn25Text
        }
    } elseif ($txt.EndsWith("a")) {
        # Drop the stray trailing "a" typo.
        $fn.Range.Text = $fn29Text
    }
}
